$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("savedEachDay")
$ws.Activate()

# ---------------------------------------------------------------------------
# Insert 10 new rows right before the existing row 291 (which holds
# "cWaterFactorForDenit"), shifting it (and everything that followed it,
# nothing in this sheet) down to row 301.
# ---------------------------------------------------------------------------
$ws.Rows.Item(291).Resize(10).Insert()

# Copy the formatting (styles) of row 290 ("cActualTranspirableWater.10")
# into the 10 freshly inserted rows, so the new entries look like the rest
# of that block.
$fmtSrc = $ws.Range("A290:I290")
$fmtDst = $ws.Range("A291:I300")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Column A on row 290 carries a bottom+left border (it is the last row of
# the previous block); the new rows should not have that border.
$ws.Range("A291:A300").Borders.LineStyle = -4142 # xlLineStyleNone

# ---------------------------------------------------------------------------
# Fill in the data for the 10 new "cFractionTranspirableWater.N" rows.
# Values are assigned in the same order the original author entered them so
# that new shared-string entries land in the same order inside
# xl/sharedStrings.xml.
# ---------------------------------------------------------------------------

# --- row 291 ---------------------------------------------------------------
$ws.Cells.Item(291,1).Value2 = "cFractionTranspirableWater.1"
$ws.Cells.Item(291,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 11"
$ws.Cells.Item(291,6).Value2 = "Fraction transpirable soil water in layer 1"
$ws.Cells.Item(291,2).Value2 = "computed"
$ws.Cells.Item(291,3).Value2 = "numeric"
$ws.Cells.Item(291,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(291,5).Value2 = "mm"
$ws.Cells.Item(291,7).Value2 = "ATSW"
$ws.Cells.Item(291,9).Value2 = 0

# --- row 292 ---------------------------------------------------------------
$ws.Cells.Item(292,1).Value2 = "cFractionTranspirableWater.2"
$ws.Cells.Item(292,6).Value2 = "Fraction transpirable soil water in layer 2"
$ws.Cells.Item(292,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 12"
$ws.Cells.Item(292,2).Value2 = "computed"
$ws.Cells.Item(292,3).Value2 = "numeric"
$ws.Cells.Item(292,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(292,5).Value2 = "mm"
$ws.Cells.Item(292,7).Value2 = "ATSW"
$ws.Cells.Item(292,9).Value2 = 0

# --- row 293 ---------------------------------------------------------------
$ws.Cells.Item(293,1).Value2 = "cFractionTranspirableWater.3"
$ws.Cells.Item(293,6).Value2 = "Fraction transpirable soil water in layer 3"
$ws.Cells.Item(293,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 13"
$ws.Cells.Item(293,2).Value2 = "computed"
$ws.Cells.Item(293,3).Value2 = "numeric"
$ws.Cells.Item(293,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(293,5).Value2 = "mm"
$ws.Cells.Item(293,7).Value2 = "ATSW"
$ws.Cells.Item(293,9).Value2 = 0

# --- row 294 ---------------------------------------------------------------
$ws.Cells.Item(294,1).Value2 = "cFractionTranspirableWater.4"
$ws.Cells.Item(294,6).Value2 = "Fraction transpirable soil water in layer 4"
$ws.Cells.Item(294,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 14"
$ws.Cells.Item(294,2).Value2 = "computed"
$ws.Cells.Item(294,3).Value2 = "numeric"
$ws.Cells.Item(294,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(294,5).Value2 = "mm"
$ws.Cells.Item(294,7).Value2 = "ATSW"
$ws.Cells.Item(294,9).Value2 = 0

# --- row 295 ---------------------------------------------------------------
$ws.Cells.Item(295,1).Value2 = "cFractionTranspirableWater.5"
$ws.Cells.Item(295,6).Value2 = "Fraction transpirable soil water in layer 5"
$ws.Cells.Item(295,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 15"
$ws.Cells.Item(295,2).Value2 = "computed"
$ws.Cells.Item(295,3).Value2 = "numeric"
$ws.Cells.Item(295,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(295,5).Value2 = "mm"
$ws.Cells.Item(295,7).Value2 = "ATSW"
$ws.Cells.Item(295,9).Value2 = 0

# --- row 296 ---------------------------------------------------------------
$ws.Cells.Item(296,1).Value2 = "cFractionTranspirableWater.6"
$ws.Cells.Item(296,6).Value2 = "Fraction transpirable soil water in layer 6"
$ws.Cells.Item(296,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 16"
$ws.Cells.Item(296,2).Value2 = "computed"
$ws.Cells.Item(296,3).Value2 = "numeric"
$ws.Cells.Item(296,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(296,5).Value2 = "mm"
$ws.Cells.Item(296,7).Value2 = "ATSW"
$ws.Cells.Item(296,9).Value2 = 0

# --- row 297 ---------------------------------------------------------------
$ws.Cells.Item(297,1).Value2 = "cFractionTranspirableWater.7"
$ws.Cells.Item(297,6).Value2 = "Fraction transpirable soil water in layer 7"
$ws.Cells.Item(297,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 17"
$ws.Cells.Item(297,2).Value2 = "computed"
$ws.Cells.Item(297,3).Value2 = "numeric"
$ws.Cells.Item(297,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(297,5).Value2 = "mm"
$ws.Cells.Item(297,7).Value2 = "ATSW"
$ws.Cells.Item(297,9).Value2 = 0

# --- row 298 ---------------------------------------------------------------
$ws.Cells.Item(298,1).Value2 = "cFractionTranspirableWater.8"
$ws.Cells.Item(298,6).Value2 = "Fraction transpirable soil water in layer 8"
$ws.Cells.Item(298,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 18"
$ws.Cells.Item(298,2).Value2 = "computed"
$ws.Cells.Item(298,3).Value2 = "numeric"
$ws.Cells.Item(298,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(298,5).Value2 = "mm"
$ws.Cells.Item(298,7).Value2 = "ATSW"
$ws.Cells.Item(298,9).Value2 = 0

# --- row 299 ---------------------------------------------------------------
$ws.Cells.Item(299,1).Value2 = "cFractionTranspirableWater.9"
$ws.Cells.Item(299,6).Value2 = "Fraction transpirable soil water in layer 9"
$ws.Cells.Item(299,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 19"
$ws.Cells.Item(299,2).Value2 = "computed"
$ws.Cells.Item(299,3).Value2 = "numeric"
$ws.Cells.Item(299,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(299,5).Value2 = "mm"
$ws.Cells.Item(299,7).Value2 = "ATSW"
$ws.Cells.Item(299,9).Value2 = 0

# --- row 300 ---------------------------------------------------------------
$ws.Cells.Item(300,1).Value2 = "cFractionTranspirableWater.10"
$ws.Cells.Item(300,6).Value2 = "Fraction transpirable soil water in layer 10"
$ws.Cells.Item(300,8).Value2 = "Quantité effective d'eau accessible pour absorption dans la couche 20"
$ws.Cells.Item(300,2).Value2 = "computed"
$ws.Cells.Item(300,3).Value2 = "numeric"
$ws.Cells.Item(300,4).Value2 = "rUpdateSoilNitrogen/rUpdateWaterBudget"
$ws.Cells.Item(300,5).Value2 = "mm"
$ws.Cells.Item(300,7).Value2 = "ATSW"
$ws.Cells.Item(300,9).Value2 = 0

# ---------------------------------------------------------------------------
# Reflect where the author ended up looking at the sheet after the edit:
# zoomed in a bit more and selected cell F300 (bottom of the new block).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("F300").Select()
